# Summary.xlsx edit: introduce a new leading scenario sheet "0" that is a
# truncated clone of sheet "0.1" (columns AP:AY of rows 3-5 removed, plus a
# handful of re-evaluated NPV figures), keep "0.1".."0.4" as-is, and drop
# the trailing "0.5" / "0.6" scenario sheets (workbook only needs 5 tabs now).

$wb = $excel.ActiveWorkbook

# --- 1. Clone "0.1" and place the clone in front of it -----------------
$source = $wb.Worksheets.Item("0.1")
$source.Copy($source)               # new copy is inserted immediately before $source
$newSheet = $wb.Worksheets.Item(1)  # the clone is now the first sheet
$newSheet.Name = "0"

# --- 2. Truncate the clone's rows 3-5 to columns B:AO -------------------
# (drop AP:AY entirely so the cells become blank, not zero)
$newSheet.Range("AP3:AY5").ClearContents()

# --- 3. Re-evaluated NPV values that differ slightly from the source ----
$newSheet.Range("AL5").Value = 3005305.231344441
$newSheet.Range("AM5").Value = 2896165.291109231
$newSheet.Range("AN5").Value = 2949795.88176504
$newSheet.Range("AO5").Value = 3005318.699777436

# --- 4. Drop the two trailing scenario sheets ---------------------------
$wb.Worksheets.Item("0.5").Delete()
$wb.Worksheets.Item("0.6").Delete()

# --- 5. Keep the first sheet ("0") as the active tab, matching the source
$newSheet.Activate()
